# Generate Report for Handoff
# Insert a new tracked file "5656ce43-a82f-4393-9d03-229b22404b63.md" as a
# new row (between "ee8ee80d-...md" and "8a9aaca5-...md") across all three
# worksheets (Overview, zh-cn, de-de), with status "Ready for handoff" /
# the associated handoff-file metadata.

$wb = $excel.ActiveWorkbook

$newMd    = "5656ce43-a82f-4393-9d03-229b22404b63.md"
$newXlfZh = "5656ce43-a82f-4393-9d03-229b22404b63.d92fc00b4dd58480f24d4d69ae4b0f873eb534b9.zh-cn.xlf"
$newXlfDe = "5656ce43-a82f-4393-9d03-229b22404b63.d92fc00b4dd58480f24d4d69ae4b0f873eb534b9.de-de.xlf"

$newMdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/8f5d986dcdcbc53d5950fb1555b1f281f71337c0/e2e/5656ce43-a82f-4393-9d03-229b22404b63.md"
$newXlfZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4b31c23766e2b34fa2a5124c78b146e5c20c5ba2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/5656ce43-a82f-4393-9d03-229b22404b63.d92fc00b4dd58480f24d4d69ae4b0f873eb534b9.zh-cn.xlf"
$newXlfDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/01679e75edf61214a860fd1e9300a5e767f6d766/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/5656ce43-a82f-4393-9d03-229b22404b63.d92fc00b4dd58480f24d4d69ae4b0f873eb534b9.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet "Overview": 3 columns (File Name / zh-cn / de-de), new row 6
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$wsOv.Rows.Item(6).Insert()

$wsOv.Range("A6").Value = $newMd
$wsOv.Range("B6").Value = "Ready for handoff"
$wsOv.Range("C6").Value = "Ready for handoff"

# Rebuild all hyperlinks on this sheet top-to-bottom so relationship ids
# come out in the same order Excel would assign them.
$wsOv.Cells.Hyperlinks.Delete()
$wsOv.Hyperlinks.Add($wsOv.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fdf01e4859064d4e6151924b2e22752896aca50e/e2e/068c0ebe-0677-4a14-a17b-8b7bc151d3c1.md", "", "", "068c0ebe-0677-4a14-a17b-8b7bc151d3c1.md")
$wsOv.Hyperlinks.Add($wsOv.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/df95368b9b7f7f52704cf10fa426e2bd97fd0f62/e2e/06aa7541-cd06-465c-8316-7632d9c3aa5a.md", "", "", "06aa7541-cd06-465c-8316-7632d9c3aa5a.md")
$wsOv.Hyperlinks.Add($wsOv.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8198ed200cbf95d4cb5860e8861c651fdabbb379/e2e/8485c808-1402-472f-9bea-7f332e27c267.md", "", "", "8485c808-1402-472f-9bea-7f332e27c267.md")
$wsOv.Hyperlinks.Add($wsOv.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/df95368b9b7f7f52704cf10fa426e2bd97fd0f62/e2e/ee8ee80d-703e-4c0c-bbc1-915a35bae61a.md", "", "", "ee8ee80d-703e-4c0c-bbc1-915a35bae61a.md")
$wsOv.Hyperlinks.Add($wsOv.Range("A6"), $newMdUrl, "", "", $newMd)
$wsOv.Hyperlinks.Add($wsOv.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/c7cf1f21ded2d375254358b64e02498a49dd54cd/e2e/8a9aaca5-1a5a-41e5-80fa-ebad684a9799.md", "", "", "8a9aaca5-1a5a-41e5-80fa-ebad684a9799.md")
$wsOv.Hyperlinks.Add($wsOv.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/ecbf7159c01a8a1b092c417d7713ddb371578928/e2e/df5dc8fc-6503-4f5c-bda6-34aa06177b80.md", "", "", "df5dc8fc-6503-4f5c-bda6-34aa06177b80.md")
$wsOv.Hyperlinks.Add($wsOv.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/ecbf7159c01a8a1b092c417d7713ddb371578928/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "zh-cn": 9 columns (detail report), new row 6
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Rows.Item(6).Insert()

$wsZh.Range("A6").Value = $newMd
$wsZh.Range("B6").Value = "Ready for handoff"
$wsZh.Range("C6").Value = $newXlfZh
$wsZh.Range("D6").Value = "2016-02-17 09:27:38"
$wsZh.Range("G6").Value = "0001-01-01 00:00:00"
$wsZh.Range("H6").Value = "Include"

$wsZh.Cells.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fdf01e4859064d4e6151924b2e22752896aca50e/e2e/068c0ebe-0677-4a14-a17b-8b7bc151d3c1.md", "", "", "068c0ebe-0677-4a14-a17b-8b7bc151d3c1.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7104a519971d0619a559f345025852723bd31c75/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/068c0ebe-0677-4a14-a17b-8b7bc151d3c1.995173ed74fdc9567a7fc9d49c2f45c1e401d036.zh-cn.xlf", "", "", "068c0ebe-0677-4a14-a17b-8b7bc151d3c1.995173ed74fdc9567a7fc9d49c2f45c1e401d036.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/df95368b9b7f7f52704cf10fa426e2bd97fd0f62/e2e/06aa7541-cd06-465c-8316-7632d9c3aa5a.md", "", "", "06aa7541-cd06-465c-8316-7632d9c3aa5a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d2a4f9b09c580b60e53a46a87f7b1b82a6bd74fc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/06aa7541-cd06-465c-8316-7632d9c3aa5a.4b62372e55ad77993176ec931bb877cbde5f71a0.zh-cn.xlf", "", "", "06aa7541-cd06-465c-8316-7632d9c3aa5a.4b62372e55ad77993176ec931bb877cbde5f71a0.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8198ed200cbf95d4cb5860e8861c651fdabbb379/e2e/8485c808-1402-472f-9bea-7f332e27c267.md", "", "", "8485c808-1402-472f-9bea-7f332e27c267.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/558c5cfc0978a6a09ef3fa1fe995b4f887f2db62/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/8485c808-1402-472f-9bea-7f332e27c267.b3febd413ba96068ef9e06bd2f1f05bf60722bd7.zh-cn.xlf", "", "", "8485c808-1402-472f-9bea-7f332e27c267.b3febd413ba96068ef9e06bd2f1f05bf60722bd7.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/47018231f7b8fe42cd09eeda4d741e47d8edb543/e2e/8485c808-1402-472f-9bea-7f332e27c267.md", "", "", "8485c808-1402-472f-9bea-7f332e27c267.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/157cb60d0c117804618cd74906de37c3acf707b9/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/8485c808-1402-472f-9bea-7f332e27c267.b3febd413ba96068ef9e06bd2f1f05bf60722bd7.zh-cn.xlf", "", "", "8485c808-1402-472f-9bea-7f332e27c267.b3febd413ba96068ef9e06bd2f1f05bf60722bd7.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/df95368b9b7f7f52704cf10fa426e2bd97fd0f62/e2e/ee8ee80d-703e-4c0c-bbc1-915a35bae61a.md", "", "", "ee8ee80d-703e-4c0c-bbc1-915a35bae61a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d2a4f9b09c580b60e53a46a87f7b1b82a6bd74fc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ee8ee80d-703e-4c0c-bbc1-915a35bae61a.05b919cc8f0ec86e6ac3a98cea0d9a70ef5e14bf.zh-cn.xlf", "", "", "ee8ee80d-703e-4c0c-bbc1-915a35bae61a.05b919cc8f0ec86e6ac3a98cea0d9a70ef5e14bf.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A6"), $newMdUrl, "", "", $newMd)
$wsZh.Hyperlinks.Add($wsZh.Range("C6"), $newXlfZhUrl, "", "", $newXlfZh)
$wsZh.Hyperlinks.Add($wsZh.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/c7cf1f21ded2d375254358b64e02498a49dd54cd/e2e/8a9aaca5-1a5a-41e5-80fa-ebad684a9799.md", "", "", "8a9aaca5-1a5a-41e5-80fa-ebad684a9799.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aad38a6a000b36672fbe87549cee9edf619f2db1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/8a9aaca5-1a5a-41e5-80fa-ebad684a9799.6c2925e9deb4fffac6eb59c95d2ddf7801228231.zh-cn.xlf", "", "", "8a9aaca5-1a5a-41e5-80fa-ebad684a9799.6c2925e9deb4fffac6eb59c95d2ddf7801228231.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/ecbf7159c01a8a1b092c417d7713ddb371578928/e2e/df5dc8fc-6503-4f5c-bda6-34aa06177b80.md", "", "", "df5dc8fc-6503-4f5c-bda6-34aa06177b80.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/38a1b102e2fb256bbdec9698c6e0b4c48ec5be95/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/df5dc8fc-6503-4f5c-bda6-34aa06177b80.592f905c61dab8d98c0fcf31a56b9e59ab62f35c.zh-cn.xlf", "", "", "df5dc8fc-6503-4f5c-bda6-34aa06177b80.592f905c61dab8d98c0fcf31a56b9e59ab62f35c.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/ecbf7159c01a8a1b092c417d7713ddb371578928/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "de-de": 9 columns (detail report), new row 6
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Rows.Item(6).Insert()

$wsDe.Range("A6").Value = $newMd
$wsDe.Range("B6").Value = "Ready for handoff"
$wsDe.Range("C6").Value = $newXlfDe
$wsDe.Range("D6").Value = "2016-02-17 09:27:50"
$wsDe.Range("G6").Value = "0001-01-01 00:00:00"
$wsDe.Range("H6").Value = "Include"

$wsDe.Cells.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fdf01e4859064d4e6151924b2e22752896aca50e/e2e/068c0ebe-0677-4a14-a17b-8b7bc151d3c1.md", "", "", "068c0ebe-0677-4a14-a17b-8b7bc151d3c1.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/127dd770976fff54e56351714d934f06aa432f9c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/068c0ebe-0677-4a14-a17b-8b7bc151d3c1.995173ed74fdc9567a7fc9d49c2f45c1e401d036.de-de.xlf", "", "", "068c0ebe-0677-4a14-a17b-8b7bc151d3c1.995173ed74fdc9567a7fc9d49c2f45c1e401d036.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/df95368b9b7f7f52704cf10fa426e2bd97fd0f62/e2e/06aa7541-cd06-465c-8316-7632d9c3aa5a.md", "", "", "06aa7541-cd06-465c-8316-7632d9c3aa5a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2cd08869fc291d43ec83138b8962b7c569e7e84b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/06aa7541-cd06-465c-8316-7632d9c3aa5a.4b62372e55ad77993176ec931bb877cbde5f71a0.de-de.xlf", "", "", "06aa7541-cd06-465c-8316-7632d9c3aa5a.4b62372e55ad77993176ec931bb877cbde5f71a0.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8198ed200cbf95d4cb5860e8861c651fdabbb379/e2e/8485c808-1402-472f-9bea-7f332e27c267.md", "", "", "8485c808-1402-472f-9bea-7f332e27c267.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/47cd4b0ba0efc28479300d571da3a1fc4368f04f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/8485c808-1402-472f-9bea-7f332e27c267.b3febd413ba96068ef9e06bd2f1f05bf60722bd7.de-de.xlf", "", "", "8485c808-1402-472f-9bea-7f332e27c267.b3febd413ba96068ef9e06bd2f1f05bf60722bd7.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5916888e8fc5d119a22c110aef5fd542683b20bb/e2e/8485c808-1402-472f-9bea-7f332e27c267.md", "", "", "8485c808-1402-472f-9bea-7f332e27c267.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/cba36bc29fedf0916313c26ca1221a9caa9d61a0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/8485c808-1402-472f-9bea-7f332e27c267.b3febd413ba96068ef9e06bd2f1f05bf60722bd7.de-de.xlf", "", "", "8485c808-1402-472f-9bea-7f332e27c267.b3febd413ba96068ef9e06bd2f1f05bf60722bd7.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/df95368b9b7f7f52704cf10fa426e2bd97fd0f62/e2e/ee8ee80d-703e-4c0c-bbc1-915a35bae61a.md", "", "", "ee8ee80d-703e-4c0c-bbc1-915a35bae61a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2cd08869fc291d43ec83138b8962b7c569e7e84b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ee8ee80d-703e-4c0c-bbc1-915a35bae61a.05b919cc8f0ec86e6ac3a98cea0d9a70ef5e14bf.de-de.xlf", "", "", "ee8ee80d-703e-4c0c-bbc1-915a35bae61a.05b919cc8f0ec86e6ac3a98cea0d9a70ef5e14bf.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A6"), $newMdUrl, "", "", $newMd)
$wsDe.Hyperlinks.Add($wsDe.Range("C6"), $newXlfDeUrl, "", "", $newXlfDe)
$wsDe.Hyperlinks.Add($wsDe.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/c7cf1f21ded2d375254358b64e02498a49dd54cd/e2e/8a9aaca5-1a5a-41e5-80fa-ebad684a9799.md", "", "", "8a9aaca5-1a5a-41e5-80fa-ebad684a9799.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6865f0c581f7281ddf1e61eaaa20aed7fb4c29ef/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/8a9aaca5-1a5a-41e5-80fa-ebad684a9799.6c2925e9deb4fffac6eb59c95d2ddf7801228231.de-de.xlf", "", "", "8a9aaca5-1a5a-41e5-80fa-ebad684a9799.6c2925e9deb4fffac6eb59c95d2ddf7801228231.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/ecbf7159c01a8a1b092c417d7713ddb371578928/e2e/df5dc8fc-6503-4f5c-bda6-34aa06177b80.md", "", "", "df5dc8fc-6503-4f5c-bda6-34aa06177b80.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b80de3a212ecb5aab645427ad7a151991371a2f1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/df5dc8fc-6503-4f5c-bda6-34aa06177b80.592f905c61dab8d98c0fcf31a56b9e59ab62f35c.de-de.xlf", "", "", "df5dc8fc-6503-4f5c-bda6-34aa06177b80.592f905c61dab8d98c0fcf31a56b9e59ab62f35c.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/ecbf7159c01a8a1b092c417d7713ddb371578928/.localization-config", "", "", ".localization-config")
